# LOB1018.xlsx edit
#
# The course-description sheet had several rows of content out of sync with
# their row labels (col A) - e.g. the "Objetivos:" row held the professor's
# name instead of the actual objectives text, "Programa:" held a date, etc.
# This edit straightens everything out: it inserts the missing
# "Objetivos:"/"Objectives:" paragraphs, shifts "Docentes responsáveis:",
# "Programa resumido:", "Programa:" and friends down to their correct rows,
# and appends a new "Bibliografia:" row (22) with the reference list that
# was missing before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 1-9: unchanged content, rewritten defensively for idempotency ---
$ws.Range("B1").Value = 'Ementa atual:'
$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'
$ws.Range("B2").Value = 'LOB1018'
$ws.Range("C2").Value = 'LOB1018'
$ws.Range("A3").Value = 'Nome:'
$ws.Range("B3").Value = ' Física I'
$ws.Range("C3").Value = ' Física I'
$ws.Range("A4").Value = 'Name:'
$ws.Range("B4").Value = 'Physics I'
$ws.Range("C4").Value = 'Physics I'
$ws.Range("A5").Value = 'Créditos-aula:'
$ws.Range("B5").Value = '4'
$ws.Range("C5").Value = '4'
$ws.Range("A6").Value = 'Créditos-trabalho'
$ws.Range("B6").Value = '0'
$ws.Range("C6").Value = '0'
$ws.Range("A7").Value = 'Carga horária:'
$ws.Range("B7").Value = '60 h'
$ws.Range("C7").Value = '60 h'
$ws.Range("A8").Value = 'Ativação:'
$ws.Range("B8").Value = '01/01/2018'
$ws.Range("C8").Value = '01/01/2018'
$ws.Range("A9").Value = 'Semestre ideal:'
$ws.Range("B9").Value = 'EF-1,EM-1,EA-2,EB-2,EP-2,EQD-1,EQN-2'
$ws.Range("C9").Value = 'EF-1,EM-1,EA-2,EB-2,EP-2,EQD-1,EQN-2'

# --- row 10: "Objetivos:" now carries the real (Portuguese) objectives text
#     instead of the professor name that had been pasted in by mistake ---
$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = 'Apresentar aos alunos os conceitos introdutórios de Física e em particular, da Mecânica incluindo cinemática e dinâmica, além de conceitos de estatística básica e análise de dados.'
$ws.Range("C10").Value = 'Apresentar aos alunos os conceitos introdutórios de Física e em particular, da Mecânica incluindo cinemática e dinâmica, além de conceitos de estatística básica e análise de dados.'

# --- row 11: "Objectives:" (English) - unchanged ---
$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").Value = 'Presenting to the students the introductory concepts of Physics and in particular, of Mechanics including kinematics and dynamics, including basic concepts of statistical and data analysis.'
$ws.Range("C11").Value = 'Presenting to the students the introductory concepts of Physics and in particular, of Mechanics including kinematics and dynamics, including basic concepts of statistical and data analysis.'

# --- row 12: "Docentes responsáveis:" label only (B/C stay empty) ---
$ws.Range("A12").Value = 'Docentes responsáveis:'

# --- row 13: professor name now sits under "Docentes responsáveis:" with no
#     label of its own in column A ---
$ws.Range("A13").Clear()
$ws.Range("B13").Value = '8711686 - Flavia Reis Cardoso Rojas'
$ws.Range("C13").Value = '8711686 - Flavia Reis Cardoso Rojas'

# --- row 14: "Programa resumido:" now correctly holds the short syllabus ---
$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = 'Introdução a física, Cinemática, Dinâmica, Trabalho, Torque e Momento Angular.'
$ws.Range("C14").Value = 'Introdução a física, Cinemática, Dinâmica, Trabalho, Torque e Momento Angular.'

# --- row 15: "Short syllabus:" (English) - same short-syllabus text ---
$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("B15").Value = 'Introdução a física, Cinemática, Dinâmica, Trabalho, Torque e Momento Angular.'
$ws.Range("C15").Value = 'Introdução a física, Cinemática, Dinâmica, Trabalho, Torque e Momento Angular.'

# --- row 16: "Programa:" now holds the full (Portuguese) syllabus instead of
#     a stray date ---
$ws.Range("A16").Value = 'Programa:'
$ws.Range("B16").Value = '1) Introdução a Física: noções de algarismos, análise dimensional, sistemas de unidades.2) Cinemática: leis de Newton e aplicações.3) Trabalho: conservação de energia, forças conservativas, aplicações.4) Impulso: momento linear e conservação de momento linear.5) Torque e momento angular: conservação de momento angular, pêndulo.'
$ws.Range("C16").Value = '1) Introdução a Física: noções de algarismos, análise dimensional, sistemas de unidades.2) Cinemática: leis de Newton e aplicações.3) Trabalho: conservação de energia, forças conservativas, aplicações.4) Impulso: momento linear e conservação de momento linear.5) Torque e momento angular: conservação de momento angular, pêndulo.'

# --- row 17 (new): "Syllabus:" (English) with the full English syllabus ---
$ws.Range("A17").Value = 'Syllabus:'
$ws.Range("B17").Value = '1) Introduction to Physics:  significant algharisms, dimensional analysis, units systems.2) Kinematics: Newton''s laws and applications.3) Work: energy conservation, conservative forces, applications.4) Impulse: momentum and conservation.5) Torque and Angular Momentum: angular momentum conservation, pendulum.'
$ws.Range("C17").Value = '1) Introduction to Physics:  significant algharisms, dimensional analysis, units systems.2) Kinematics: Newton''s laws and applications.3) Work: energy conservation, conservative forces, applications.4) Impulse: momentum and conservation.5) Torque and Angular Momentum: angular momentum conservation, pendulum.'

# --- row 18: "Avaliação:" label only (B/C now cleared - they mistakenly held
#     the professor's name) ---
$ws.Range("A18").Value = 'Avaliação:'
$ws.Range("B18:C18").Clear()

# --- row 19: "Método:" now holds the evaluation-method paragraph ---
$ws.Range("A19").Value = 'Método:'
$ws.Range("B19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range("C19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'

# --- row 20: "Critério:" now holds the passing-grade criterion ---
$ws.Range("A20").Value = 'Critério:'
$ws.Range("B20").Value = 'NF≥ 5,0.'
$ws.Range("C20").Value = 'NF≥ 5,0.'

# --- row 21: "Norma de recuperação:" now holds the make-up-exam rule ---
$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Range("B21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada'
$ws.Range("C21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada'

# --- row 22 (new): "Bibliografia:" with the reading list ---
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = 'NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 1, Edgard Blucher (2008).RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol.1, LTC (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.1, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 1, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 1, Thomson Pioneira (2008).'
$ws.Range("C22").Value = 'NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 1, Edgard Blucher (2008).RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol.1, LTC (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.1, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 1, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 1, Thomson Pioneira (2008).'

# --- fix up row heights: rows 13 & 18 go back to the default (no label /
#     no long text anymore), rows 15 & 21 shrink to 60pt, rows 17 & 22 pick
#     up the 60pt/120pt heights matching their new long-text content ---
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
